# Update "Inscritos" (E), "Pagos" (F) and "Inscrições homologadas" (H) counts
# for a handful of rows in the Inscricoes sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("E2").Value = 40

$ws.Range("E8").Value = 10

$ws.Range("E15").Value = 138

$ws.Range("E24").Value = 20
$ws.Range("F24").Value = 11
$ws.Range("H24").Value = 11

$ws.Range("E36").Value = 75

$ws.Range("E37").Value = 41

$ws.Range("E41").Value = 28

$ws.Range("E42").Value = 29

$ws.Range("E50").Value = 20
$ws.Range("F50").Value = 4
$ws.Range("H50").Value = 4

$ws.Range("E67").Value = 34

$ws.Range("E70").Value = 32

$ws.Range("E73").Value = 22
$ws.Range("F73").Value = 7
$ws.Range("H73").Value = 7

$wb.Save()
